$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.648.41'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.083.17'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.79%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.55%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.41'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.99'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +7.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.077.27'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.51%  '

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.03%  '

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.29'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.29'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +7.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.566.80'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.596.47'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.067.42'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.109'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.77'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '487.99'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.54'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.24'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +7.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.60'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.41'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +6.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.76'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.02'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.63%  '

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.99'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +10.19%  '

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.20'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.99%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.46%  '

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.81'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +10.83%  '

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Stacks'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.43'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +9.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.67'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.00'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '469.65'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.192.37'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.27%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0818'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0398'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.121'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.26'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.55'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.82'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +13.08%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.253'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.05%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.05'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +8.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.110'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0525'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.99%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '116.34'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.16%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.08'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.36%  '
